$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "Tue Jan 28 21:57:43 EST 2025"
$ws.Cells.Item(3, 2).Value = "Tue Jan 28 21:57:55 EST 2025"
$ws.Cells.Item(4, 2).Value = "Tue Jan 28 21:58:07 EST 2025"
$ws.Cells.Item(5, 2).Value = "Tue Jan 28 21:58:19 EST 2025"
$ws.Cells.Item(6, 2).Value = "Tue Jan 28 21:58:30 EST 2025"
$ws.Cells.Item(7, 2).Value = "Tue Jan 28 21:58:42 EST 2025"
$ws.Cells.Item(8, 2).Value = "Tue Jan 28 21:58:54 EST 2025"
$ws.Cells.Item(9, 2).Value = "Tue Jan 28 21:59:06 EST 2025"
$ws.Cells.Item(10, 2).Value = "Tue Jan 28 21:59:18 EST 2025"
$ws.Cells.Item(11, 2).Value = "Tue Jan 28 21:59:30 EST 2025"
$ws.Cells.Item(12, 2).Value = "Tue Jan 28 21:59:42 EST 2025"
$ws.Cells.Item(13, 2).Value = "Tue Jan 28 21:59:54 EST 2025"
$ws.Cells.Item(14, 2).Value = "Tue Jan 28 22:00:06 EST 2025"
$ws.Cells.Item(15, 2).Value = "Tue Jan 28 22:00:18 EST 2025"
$ws.Cells.Item(16, 2).Value = "Tue Jan 28 22:00:31 EST 2025"
$ws.Cells.Item(17, 2).Value = "Tue Jan 28 22:00:44 EST 2025"
$ws.Cells.Item(18, 2).Value = "Tue Jan 28 22:00:57 EST 2025"
$ws.Cells.Item(19, 2).Value = "Tue Jan 28 22:01:10 EST 2025"
$ws.Cells.Item(20, 2).Value = "Tue Jan 28 22:01:23 EST 2025"
$ws.Cells.Item(21, 2).Value = "Tue Jan 28 22:01:36 EST 2025"
$ws.Cells.Item(22, 2).Value = "Tue Jan 28 22:01:48 EST 2025"
$ws.Cells.Item(23, 2).Value = "Tue Jan 28 22:02:01 EST 2025"
$ws.Cells.Item(24, 2).Value = "Tue Jan 28 22:02:14 EST 2025"
$ws.Cells.Item(25, 2).Value = "Tue Jan 28 22:02:26 EST 2025"
$ws.Cells.Item(26, 2).Value = "Tue Jan 28 22:02:39 EST 2025"
$ws.Cells.Item(27, 2).Value = "Tue Jan 28 22:02:51 EST 2025"
$ws.Cells.Item(28, 2).Value = "Tue Jan 28 22:03:04 EST 2025"
$ws.Cells.Item(29, 2).Value = "Tue Jan 28 22:03:16 EST 2025"
$ws.Cells.Item(30, 2).Value = "Tue Jan 28 22:03:29 EST 2025"
$ws.Cells.Item(31, 2).Value = "Tue Jan 28 22:03:42 EST 2025"
$ws.Cells.Item(32, 2).Value = "Tue Jan 28 22:03:56 EST 2025"
$ws.Cells.Item(33, 2).Value = "Tue Jan 28 22:04:09 EST 2025"
$ws.Cells.Item(34, 2).Value = "Tue Jan 28 22:04:22 EST 2025"
$ws.Cells.Item(36, 2).Value = "Tue Jan 28 22:04:34 EST 2025"
$ws.Cells.Item(37, 2).Value = "Tue Jan 28 22:04:46 EST 2025"
$ws.Cells.Item(38, 2).Value = "Tue Jan 28 22:04:58 EST 2025"
$ws.Cells.Item(39, 2).Value = "Tue Jan 28 22:05:10 EST 2025"
$ws.Cells.Item(40, 2).Value = "Tue Jan 28 22:05:22 EST 2025"
$ws.Cells.Item(41, 2).Value = "Tue Jan 28 22:05:34 EST 2025"
$ws.Cells.Item(42, 2).Value = "Tue Jan 28 22:05:46 EST 2025"
$ws.Cells.Item(43, 2).Value = "Tue Jan 28 22:05:59 EST 2025"
$ws.Cells.Item(44, 2).Value = "Tue Jan 28 22:06:12 EST 2025"
$ws.Cells.Item(45, 2).Value = "Tue Jan 28 22:06:25 EST 2025"
$ws.Cells.Item(47, 2).Value = "Tue Jan 28 22:06:38 EST 2025"
$ws.Cells.Item(48, 2).Value = "Tue Jan 28 22:06:51 EST 2025"
$ws.Cells.Item(49, 2).Value = "Tue Jan 28 22:07:04 EST 2025"
$ws.Cells.Item(50, 2).Value = "Tue Jan 28 22:07:17 EST 2025"
$ws.Cells.Item(51, 2).Value = "Tue Jan 28 22:07:30 EST 2025"
$ws.Cells.Item(52, 2).Value = "Tue Jan 28 22:07:43 EST 2025"
$ws.Cells.Item(53, 2).Value = "Tue Jan 28 22:07:56 EST 2025"
$ws.Cells.Item(54, 2).Value = "Tue Jan 28 22:08:09 EST 2025"
